$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.9999999896980597,
    0.9979140901101089,
    0.99500911094715372,
    0.99829994061194649,
    0.98844870547848007,
    0.96720736925804229,
    0.96430219109866777,
    0.95893268788060304,
    0.95450479057813964,
    0.9509928921930928,
    0.9504154328161335,
    0.94957274327561958,
    0.94902966015798262,
    0.95016936533334118,
    0.95217291944375271,
    0.95490657617021701,
    0.95119909219602317,
    0.95009022088517103,
    0.99686094760782917,
    0.98974422730827283,
    0.98834579101362052,
    0.98708130049921083,
    0.98862830986042249,
    0.97560830668469378,
    0.96915147150716141,
    0.96460109901984914,
    0.96121625791350607,
    0.94929821801106118,
    0.93691729457385242,
    0.9316322516507074,
    0.9277326277516762,
    0.92739969671286127,
    0.92687970325976643
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
